$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Gilbert quest) ---
$ws.Range("C2").Value = 17
$ws.Range("E2").Value = "길버트 아저씨가 도와달라 한다. 얻어먹은 것도 있으니 일단은 노력해보자."
$ws.Range("G2").Value = "5;"
$ws.Range("H2").Value = "0;1;5;"
$ws.Range("I2").Value = "5;5;1;"

# --- Add new row 3 (suspicious man quest) ---
$ws.Range("A3").Value = "수상한 남자"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = "수상한 남자의 말을 들어주자"
$ws.Range("E3").Value = "수상한 살람이 다음 지역으로 가는 방법을 알고 있지만, 쉽사리 알려주지 않는다. `n특정 아이템을 요구하는 듯한데-."
$ws.Range("F3").Value = "5;"
$ws.Range("G3").Value = "1;"
$ws.Range("H3").Value = "0;"
$ws.Range("I3").Value = "5;"

# Row 3 formatting: wrap text on the description cell + taller row to fit
$ws.Range("E3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 49.5

# Column E width tweak (was bestFit 69.375 -> custom 63.25)
$ws.Columns.Item(5).ColumnWidth = 62.57

# Update selection to the new last cell, matching the saved view state
$ws.Range("I3").Select() | Out-Null
